# Reset coal retirements data to align with EPA 111 RIA
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRpUNL")
$wsAbout = $wb.Worksheets.Item("About")

# Update the "Capacity Retired per Unit Net Loss" values for the plant types
# that previously had 0.01 -> now 0.03 (rows 2-5, 7-8, 13-15)
$ws.Range("B2").Value = 0.03
$ws.Range("B3").Value = 0.03
$ws.Range("B4").Value = 0.03
$ws.Range("B5").Value = 0.03
$ws.Range("B7").Value = 0.03
$ws.Range("B8").Value = 0.03
$ws.Range("B13").Value = 0.03
$ws.Range("B14").Value = 0.03
$ws.Range("B15").Value = 0.03

# Update the active cell selection on the CRpUNL sheet to match the new view state
$ws.Activate()
$ws.Range("F15").Select()

# Restore the originally active sheet ("About") so the workbook's active tab
# selection is unchanged
$wsAbout.Activate()
